$wb = $excel.ActiveWorkbook

# Replace the "Ready for handoff" status text with "In Translation"
# across every sheet (Overview, zh-cn, de-de).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value()
            if ("Ready for handoff" -eq $val) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# Shrink the now-narrower "Status" columns to match the new content width.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 13.4101845877511
$overview.Columns.Item(6).ColumnWidth = 13.4101845877511

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 13.4101845877511

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 13.4101845877511
